$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix formatting of "Importe" (amount) column H: values were scraped with
# Spanish/Argentine thousands "." / decimal "," formatting (e.g. "49.600,00")
# and must become plain decimal-point numeric text (e.g. "49600.00").
# NumberFormat is forced to text ("@") before the write so Excel keeps the
# literal string instead of re-parsing it as a real number, then the cell
# style is reset back to Normal so no visible formatting changes.
$importeFixes = @{
    2 = "49600.00"
    3 = "4500.00"
    4 = "49828.00"
    5 = "33840.00"
    6 = "160065.00"
    7 = "105600.00"
    8 = "173650.00"
    9 = "872110.00"
    10 = "70860.00"
    11 = "278200.00"
    12 = "508800.00"
    13 = "59100.00"
    14 = "198008.00"
    15 = "440000.00"
    16 = "195000.00"
    17 = "190000.00"
    18 = "4489.10"
    19 = "435600.00"
    20 = "4560.00"
    21 = "484057.87"
    22 = "561768.32"
    23 = "60492.00"
    24 = "17475.00"
    25 = "256412.98"
    26 = "159440.21"
    27 = "43619.19"
    28 = "456.00"
    29 = "31855.77"
    30 = "55208.00"
    31 = "15000.00"
    32 = "30789.99"
    33 = "27645.40"
    34 = "45.00"
    35 = "1268.78"
    36 = "133.64"
    37 = "39761.00"
    38 = "146730.00"
    39 = "621100.00"
    40 = "130430.00"
    41 = "8887.84"
    42 = "7139.55"
    43 = "136.00"
    44 = "14059.64"
    45 = "7375.00"
    46 = "6350.00"
    47 = "100.44"
    48 = "5413.00"
    49 = "600.00"
    50 = "26769.00"
    51 = "1660.00"
    52 = "1894.91"
    53 = "145.00"
    54 = "80.00"
    55 = "58485.50"
    56 = "340.00"
    57 = "356896.00"
    58 = "7675.00"
    59 = "8251.87"
    60 = "1067090.00"
    61 = "860.00"
    62 = "71090.00"
    63 = "35750.00"
    64 = "94600.00"
    65 = "749100.00"
    66 = "44000.00"
    67 = "257840.00"
    68 = "206.00"
    69 = "174.55"
    70 = "2800.00"
    71 = "25480.00"
    72 = "4127.00"
    73 = "291598.00"
    74 = "3750.00"
    75 = "15074.00"
    76 = "47240.00"
    77 = "300.00"
    78 = "92070.00"
    79 = "7600.00"
    80 = "41651.02"
    81 = "18200.00"
    82 = "4790.00"
    83 = "6450.00"
    84 = "10756.35"
    85 = "28314.00"
    86 = "4050.00"
    87 = "7056.00"
    88 = "1355.00"
    89 = "2909.00"
    90 = "13499.99"
    91 = "6660.00"
    92 = "44190.00"
    93 = "64.17"
    94 = "24260.34"
    95 = "3605.90"
    96 = "14245.72"
    97 = "844.90"
    98 = "25433.20"
    99 = "38086.30"
    100 = "12198.92"
    101 = "2530.00"
    102 = "3120.00"
    103 = "4500.00"
    104 = "570.00"
    105 = "1067.88"
    106 = "1340.10"
    107 = "900.00"
    108 = "768.00"
    109 = "195.70"
    110 = "130.00"
    111 = "4000.00"
    112 = "500.00"
    113 = "9150.00"
    114 = "160000.00"
    115 = "6365.00"
    116 = "4356.00"
    117 = "439710.00"
    118 = "7850.00"
    119 = "84700.00"
    120 = "22928.00"
    121 = "3875.32"
    122 = "13400.00"
    123 = "15000.00"
    124 = "94200.00"
    125 = "35580.00"
    126 = "25000.00"
    127 = "1966250.00"
    128 = "71200.00"
    129 = "1483.97"
    130 = "4598.00"
    131 = "800.00"
    132 = "200.00"
    133 = "2040.00"
    134 = "4238.20"
    135 = "480000.00"
    136 = "7776.00"
    137 = "6080.00"
    138 = "200.00"
    139 = "15730.00"
    140 = "8687.00"
    141 = "7200.00"
    142 = "99.98"
    143 = "942.60"
    144 = "2023.68"
    145 = "7150.00"
    146 = "966.01"
    147 = "992.00"
    148 = "150.00"
    149 = "42700.00"
    150 = "15091.27"
    151 = "5125.68"
    152 = "2080.00"
    153 = "9100.00"
    154 = "24.00"
    155 = "62742.18"
    156 = "68590.00"
    157 = "26190.00"
    158 = "16272.00"
    159 = "4615.74"
    160 = "8500.00"
    161 = "30000.00"
    162 = "30000.00"
    163 = "30000.00"
    164 = "30000.00"
    165 = "30000.00"
    166 = "30000.00"
    167 = "60000.00"
    168 = "60000.00"
    169 = "60000.00"
    170 = "30000.00"
    171 = "30662.72"
    172 = "8894538.38"
    173 = "30000.00"
    174 = "381600.00"
    175 = "2300000.00"
    176 = "135500.00"
    177 = "144500.00"
    178 = "135500.00"
    179 = "135500.00"
    180 = "135500.00"
    181 = "135500.00"
    182 = "239000.00"
    183 = "245900.00"
    184 = "347000.00"
    185 = "135500.00"
    186 = "135500.00"
    187 = "135500.00"
    188 = "135500.00"
    189 = "135500.00"
    190 = "239000.00"
    191 = "342500.00"
    192 = "239000.00"
    193 = "135500.00"
    194 = "244000.00"
    195 = "135500.00"
    196 = "135500.00"
    197 = "141800.00"
    198 = "135500.00"
    199 = "80000.00"
    200 = "6620.00"
    201 = "8000.00"
    202 = "454887.50"
    203 = "100000.00"
    204 = "12705.00"
    205 = "114840.00"
    206 = "95240.00"
    207 = "92600.00"
    208 = "32400.00"
    209 = "6400.00"
}

foreach ($row in $importeFixes.Keys) {
    $cell = $ws.Cells.Item($row, 8)
    $cell.NumberFormat = "@"
    $cell.Value = $importeFixes[$row]
    $cell.Style = "Normal"
}

# Fix the "Razon social" text for row 153 (E column): commas standardised to
# periods and the "S.H." abbreviation normalised to "SH".
$ws.Range("E153").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

